# Catch_Trust_4.xlsx update:
#  - Row 2 (Arnoglossus laterna, 1-RAP) weight/number corrected.
#  - The old row 5 (Solea solea / SOLEVUL / 8.492 / 105) is removed,
#    shifting every subsequent row up by one (rows 6..55 -> 5..54),
#    shrinking the used range from A1:K55 to A1:K54.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update W(kg) and Numb for row 2 (Arnoglossus laterna, 1-RAP).
$ws.Range("G2").Value = 4.22
$ws.Range("H2").Value = 52

# Delete row 5 entirely; Excel shifts rows 6:55 up to 5:54.
$ws.Rows.Item(5).Delete()
